$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.190.94"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "2.359.73"
$ws.Range("E3").Value = "  +2.63%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").Value = "302.66"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "95.59"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").Value = "0.504"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "34.17"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "18.61"
$ws.Range("E12").Value = "  -3.43%  "
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "2.721.51"
$ws.Range("E15").Value = "  +2.65%  "
$ws.Range("D16").Value = "2.370.20"
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("D17").Value = "0.799"
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "43.167.84"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "12.21"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "6.26"
$ws.Range("E20").Value = "  +4.03%  "
$ws.Range("D21").Value = "0.0₃0891"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "68.18"
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "235.50"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("D27").Value = "24.62"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("E28").Value = "  +14.87%  "
$ws.Range("D29").Value = "9.16"
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("D30").Value = "31.35"
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D34").Value = "17.21"
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("E35").Value = "  +5.75%  "
$ws.Range("D36").Value = "4.37"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").Value = "2.32"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").Value = "22.58"
$ws.Range("E39").Value = "  +13.68%  "
$ws.Range("D40").Value = "2.77"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Value = "111.85"
$ws.Range("E42").Value = "  -32.25%  "
$ws.Range("D43").Value = "1.944.05"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").Value = "0.0282"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").Value = "2.11"
$ws.Range("E45").Value = "  +3.45%  "
$ws.Range("D46").Value = "9.42"
$ws.Range("E46").Value = "  -9.90%  "
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").Value = "2.585.51"
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("D49").Value = "52.83"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("E50").Value = "  -2.99%  "
$ws.Range("D51").Value = "72.14"
$ws.Range("E51").Value = "  +0.81%  "
